$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 4858.8765
$ws.Range("I15").Value = 4858.8765
$ws.Range("K15").Value = 14576.6295
$ws.Range("M15").Value = -14407.6295
$ws.Range("H138").Value = 5192.615
$ws.Range("I138").Value = 4500.5713
$ws.Range("J138").Value = 6000
$ws.Range("K138").Value = 13501.7139
$ws.Range("L138").Value = 18000
$ws.Range("M138").Value = -8361.713899999999
$ws.Range("N138").Value = -28280

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H31").Value = 7440.3335
$ws.Range("I31").Value = 7440.3335
$ws.Range("K31").Value = 7440.3335
$ws.Range("M31").Value = -7146.3335
$ws.Range("H45").Value = 85417.75
$ws.Range("I45").Value = 168949.83
$ws.Range("J45").Value = 1885.6666
$ws.Range("K45").Value = 168949.83
$ws.Range("L45").Value = 1885.6666
$ws.Range("M45").Value = -168572.83
$ws.Range("N45").Value = -2639.6666
$ws.Range("H61").Value = 2517.5366
$ws.Range("I61").Value = 1978.0625
$ws.Range("J61").Value = 4435.6665
$ws.Range("K61").Value = 1978.0625
$ws.Range("L61").Value = 4435.6665
$ws.Range("M61").Value = -1766.0625
$ws.Range("N61").Value = -4859.6665
$ws.Range("H62").Value = 34000
$ws.Range("J62").Value = 34000
$ws.Range("L62").Value = 34000
$ws.Range("N62").Value = -35248
$ws.Range("H65").Value = 34000
$ws.Range("J65").Value = 34000
$ws.Range("L65").Value = 102000
$ws.Range("N65").Value = -108240
$ws.Range("H102").Value = 1776
$ws.Range("I102").Value = 1907.1428
$ws.Range("J102").Value = 1408.8
$ws.Range("K102").Value = 1907.1428
$ws.Range("L102").Value = 1408.8
$ws.Range("M102").Value = -285.1428000000001
$ws.Range("N102").Value = -4652.8
$ws.Range("H132").Value = 1933.5278
$ws.Range("I132").Value = 1709.8077
$ws.Range("K132").Value = 5129.4231
$ws.Range("M132").Value = -2599.4231
$ws.Range("H136").Value = 2517.5366
$ws.Range("I136").Value = 1978.0625
$ws.Range("J136").Value = 4435.6665
$ws.Range("K136").Value = 5934.1875
$ws.Range("L136").Value = 13306.9995
$ws.Range("M136").Value = -3384.1875
$ws.Range("N136").Value = -18406.9995

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2632.8
$ws.Range("I20").Value = 2744.889
$ws.Range("J20").Value = 2464.6667
$ws.Range("K20").Value = 2744.889
$ws.Range("L20").Value = 2464.6667
$ws.Range("M20").Value = -2497.889
$ws.Range("N20").Value = -2958.6667
$ws.Range("H21").Value = 12500
$ws.Range("J21").Value = 12500
$ws.Range("L21").Value = 12500
$ws.Range("N21").Value = -12972
$ws.Range("H43").Value = 115512.75
$ws.Range("J43").Value = 115512.75
$ws.Range("L43").Value = 115512.75
$ws.Range("N43").Value = -115874.75
$ws.Range("H59").Value = 39323.637
$ws.Range("J59").Value = 39323.637
$ws.Range("L59").Value = 39323.637
$ws.Range("N59").Value = -41017.637
$ws.Range("H94").Value = 1030.5
$ws.Range("I94").Value = 893.5
$ws.Range("J94").Value = 1852.5
$ws.Range("K94").Value = 893.5
$ws.Range("L94").Value = 1852.5
$ws.Range("M94").Value = -442.5
$ws.Range("N94").Value = -2754.5
$ws.Range("H98").Value = 45000
$ws.Range("J98").Value = 45000
$ws.Range("L98").Value = 45000
$ws.Range("N98").Value = -50990
$ws.Range("H99").Value = 905
$ws.Range("I99").Value = 840
$ws.Range("J99").Value = 1100
$ws.Range("K99").Value = 840
$ws.Range("L99").Value = 1100
$ws.Range("M99").Value = 658
$ws.Range("N99").Value = -4096
$ws.Range("H107").Value = 2100.45
$ws.Range("I107").Value = 1625.5625
$ws.Range("J107").Value = 4000
$ws.Range("K107").Value = 1625.5625
$ws.Range("L107").Value = 4000
$ws.Range("M107").Value = 294.4375
$ws.Range("N107").Value = -7840
$ws.Range("H109").Value = 44875
$ws.Range("J109").Value = 44875
$ws.Range("L109").Value = 44875
$ws.Range("N109").Value = -47649
$ws.Range("H119").Value = 35000
$ws.Range("J119").Value = 35000
$ws.Range("L119").Value = 35000
$ws.Range("N119").Value = -44676

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2317.0544
$ws.Range("I58").Value = 1238.8846
$ws.Range("J58").Value = 3283.6897
$ws.Range("K58").Value = 1238.8846
$ws.Range("L58").Value = 3283.6897
$ws.Range("M58").Value = -1035.8846
$ws.Range("N58").Value = -3689.6897
$ws.Range("H97").Value = 13900
$ws.Range("J97").Value = 13900
$ws.Range("L97").Value = 13900
$ws.Range("N97").Value = -15882
$ws.Range("H132").Value = 4682.6113
$ws.Range("I132").Value = 5008.1816
$ws.Range("J132").Value = 4171
$ws.Range("K132").Value = 15024.5448
$ws.Range("L132").Value = 12513
$ws.Range("M132").Value = -12494.5448
$ws.Range("N132").Value = -17573
$ws.Range("H136").Value = 2317.0544
$ws.Range("I136").Value = 1238.8846
$ws.Range("J136").Value = 3283.6897
$ws.Range("K136").Value = 3716.6538
$ws.Range("L136").Value = 9851.069100000001
$ws.Range("M136").Value = -1166.6538
$ws.Range("N136").Value = -14951.0691
$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H133").Value = 6823.619
$ws.Range("I133").Value = 4933.3335
$ws.Range("J133").Value = 7138.6665
$ws.Range("K133").Value = 14800.0005
$ws.Range("L133").Value = 21415.9995
$ws.Range("M133").Value = -9740.000499999998
$ws.Range("N133").Value = -31535.9995

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H99").Value = 4214
$ws.Range("I99").Value = 1697.4
$ws.Range("J99").Value = 29380
$ws.Range("K99").Value = 1697.4
$ws.Range("L99").Value = 29380
$ws.Range("M99").Value = 548.5999999999999
$ws.Range("N99").Value = -33872
$ws.Range("H132").Value = 1970.8223
$ws.Range("I132").Value = 2083.2856
$ws.Range("J132").Value = 1872.4166
$ws.Range("K132").Value = 6249.8568
$ws.Range("L132").Value = 5617.2498
$ws.Range("M132").Value = -3719.8568
$ws.Range("N132").Value = -10677.2498

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2788.6667
$ws.Range("I40").Value = 2512.25
$ws.Range("K40").Value = 2512.25
$ws.Range("M40").Value = -2376.25
$ws.Range("H61").Value = 2511.5557
$ws.Range("I61").Value = 1184
$ws.Range("K61").Value = 1184
$ws.Range("M61").Value = -982
$ws.Range("H62").Value = 33000
$ws.Range("J62").Value = 33000
$ws.Range("L62").Value = 33000
$ws.Range("N62").Value = -34248
$ws.Range("H65").Value = 33000
$ws.Range("J65").Value = 33000
$ws.Range("L65").Value = 99000
$ws.Range("N65").Value = -105240
$ws.Range("H99").Value = 70000
$ws.Range("J99").Value = 70000
$ws.Range("L99").Value = 70000
$ws.Range("N99").Value = -75990
$ws.Range("H109").Value = 24500
$ws.Range("J109").Value = 24500
$ws.Range("L109").Value = 24500
$ws.Range("N109").Value = -27274
$ws.Range("H113").Value = 2511.5557
$ws.Range("I113").Value = 1184
$ws.Range("K113").Value = 1184
$ws.Range("M113").Value = 986
$ws.Range("H120").Value = 0
$ws.Range("J120").Value = 0
$ws.Range("L120").Value = 0
$ws.Range("N120").ClearContents()
$ws.Range("H130").Value = 49985.3
$ws.Range("J130").Value = 49985.3
$ws.Range("L130").Value = 49985.3
$ws.Range("N130").Value = -60025.3
$ws.Range("H133").Value = 47950
$ws.Range("J133").Value = 47950
$ws.Range("L133").Value = 47950
$ws.Range("N133").Value = -53010

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1994.1765
$ws.Range("I81").Value = 1516.8334
$ws.Range("J81").Value = 2254.5454
$ws.Range("K81").Value = 3033.6668
$ws.Range("L81").Value = 4509.0908
$ws.Range("M81").Value = -1972.6668
$ws.Range("N81").Value = -6631.0908
$ws.Range("H84").Value = 1994.1765
$ws.Range("I84").Value = 1516.8334
$ws.Range("J84").Value = 2254.5454
$ws.Range("K84").Value = 15168.334
$ws.Range("L84").Value = 22545.454
$ws.Range("M84").Value = -9864.333999999999
$ws.Range("N84").Value = -33153.454
$ws.Range("H119").Value = 30000
$ws.Range("J119").Value = 30000
$ws.Range("L119").Value = 30000
$ws.Range("N119").Value = -39676
$ws.Range("H132").Value = 5754.049
$ws.Range("I132").Value = 7533.5557
$ws.Range("J132").Value = 2322.1428
$ws.Range("K132").Value = 22600.6671
$ws.Range("L132").Value = 6966.428400000001
$ws.Range("M132").Value = -20070.6671
$ws.Range("N132").Value = -12026.4284
